# BR_Excel.xlsx update
#  - Replicating Bouyer paper: the "Cue" label in B1 is renamed to lowercase "cue".
#  - The active selection on Sheet1 is moved from C3 to D4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 currently holds the shared string "Cue" -> change to "cue"
$ws.Range("B1").Value = "cue"

# Update the sheet's active selection/cell to D4 (was C3)
$ws.Range("D4").Select() | Out-Null
